$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "author" column (B) previously held the author's initials ("mz").
# Switch every occurrence to the author's full name.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 2).Value = "Maja Založnik"
}

# Reflect the new active-cell selection recorded in the saved sheet view.
$ws.Range("B2").Select()
